$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the full target range is treated as text, matching the original
# inline-string cell formatting (so numeric-looking / date-looking values
# such as '2024-01-20' or '12.5' are not auto-converted by Excel).
$ws.Range("A8:K14").NumberFormat = "@"

# Row 8
$ws.Range("A8").Value = 'edit_check'
$ws.Range("B8").Value = 'AE'
$ws.Range("C8").Value = 'SUBJ-1008'
$ws.Range("D8").Value = 'AEENDTC'
$ws.Range("E8").Value = 'Missing end date for ongoing AE.'
$ws.Range("F8").Value = '2024-01-20'
$ws.Range("G8").Value = ''
$ws.Range("H8").Value = 'AEENDTC'
$ws.Range("I8").Value = ''
$ws.Range("J8").Value = 'Required when AEOUT=Ongoing'
$ws.Range("K8").Value = ''

# Row 9
$ws.Range("A9").Value = 'listing'
$ws.Range("B9").Value = 'LB'
$ws.Range("C9").Value = 'SUBJ-1009'
$ws.Range("D9").Value = 'LBORRES, LBCLSIG'
$ws.Range("E9").Value = 'Discrepancy vs central lab: EDC value differs from external. Clinical significance unclear.'
$ws.Range("F9").Value = ''
$ws.Range("G9").Value = ''
$ws.Range("H9").Value = 'LBORRES'
$ws.Range("I9").Value = '12.5'
$ws.Range("J9").Value = 'Central lab 11.8'
$ws.Range("K9").Value = 'LBCLSIG=Y, requires clinical context'

# Row 10
$ws.Range("A10").Value = 'edit_check'
$ws.Range("B10").Value = 'DM'
$ws.Range("C10").Value = 'SUBJ-1010'
$ws.Range("D10").Value = 'BRTHDTC'
$ws.Range("E10").Value = 'Invalid or partial date: birth date month/year only. May be acceptable per protocol.'
$ws.Range("F10").Value = ''
$ws.Range("G10").Value = ''
$ws.Range("H10").Value = 'BRTHDTC'
$ws.Range("I10").Value = '1985-03'
$ws.Range("J10").Value = 'ISO 8601 full date'
$ws.Range("K10").Value = 'Partial date may be acceptable'

# Row 11
$ws.Range("A11").Value = 'listing'
$ws.Range("B11").Value = 'AE'
$ws.Range("C11").Value = 'SUBJ-1011'
$ws.Range("D11").Value = 'AETERM, AESEV, AESER'
$ws.Range("E11").Value = 'Complex adverse event with multiple related conditions. Requires medical review to determine if single or multiple events.'
$ws.Range("F11").Value = '2024-04-10'
$ws.Range("G11").Value = ''
$ws.Range("H11").Value = 'AETERM'
$ws.Range("I11").Value = 'Headache, Nausea, Dizziness'
$ws.Range("J11").Value = 'Multiple symptoms'
$ws.Range("K11").Value = 'Need to assess if related or separate'

# Row 12
$ws.Range("A12").Value = 'listing'
$ws.Range("B12").Value = 'DM'
$ws.Range("C12").Value = 'SUBJ-1012'
$ws.Range("D12").Value = 'DMWEIGHT, DMWTU'
$ws.Range("E12").Value = 'Inconsistent units: weight in kg vs lb across visits. Need to assess impact on BMI calculations.'
$ws.Range("F12").Value = ''
$ws.Range("G12").Value = ''
$ws.Range("H12").Value = 'DMWTU'
$ws.Range("I12").Value = 'lb'
$ws.Range("J12").Value = 'Study standard kg'
$ws.Range("K12").Value = 'Previous visit was kg, BMI affected'

# Row 13
$ws.Range("A13").Value = 'edit_check'
$ws.Range("B13").Value = 'AE'
$ws.Range("C13").Value = 'SUBJ-1013'
$ws.Range("D13").Value = 'AESTDTC, AEENDTC, AESER'
$ws.Range("E13").Value = 'Serious adverse event with ambiguous timeline. Start date conflicts with hospitalization records.'
$ws.Range("F13").Value = '2024-05-01'
$ws.Range("G13").Value = '2024-05-15'
$ws.Range("H13").Value = 'AESTDTC'
$ws.Range("I13").Value = '2024-05-03'
$ws.Range("J13").Value = 'Hospitalization started 2024-05-03'
$ws.Range("K13").Value = 'Date reconciliation needed'

# Row 14
$ws.Range("A14").Value = 'listing'
$ws.Range("B14").Value = 'CM'
$ws.Range("C14").Value = 'SUBJ-1014'
$ws.Range("D14").Value = 'CMTRT, CMDOSFRM'
$ws.Range("E14").Value = 'Uncommon medication coding issue: combination product not in standard dictionary. Requires manual review.'
$ws.Range("F14").Value = ''
$ws.Range("G14").Value = ''
$ws.Range("H14").Value = 'CMTRT'
$ws.Range("I14").Value = 'Product-XY Plus'
$ws.Range("J14").Value = 'Not in MedDRA/WHODrug'
$ws.Range("K14").Value = 'Novel combination product'
